$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "OGL" row entirely (was row 4); this shifts subsequent rows up by one.
$ws.Rows.Item(4).Delete()

# Update the wxDC result value, and drop the now unused C column result.
$ws.Range("B2").Value = "17.5099999905  seconds"
$ws.Range("C2").ClearContents()

# Update the header in B1 to reflect the new test description.
$ws.Range("B1").Value = "ReDrawing 2500 Rects (w/text) 100 times"

# Clear the now-unused C1/D1 header cells (keep their bold styling, drop the text).
$ws.Range("C1:D1").ClearContents()

# Remove leftover results for FloatCanvas and Pygame (w/ wx) rows.
$ws.Range("B3:C3").ClearContents()
$ws.Range("B4:C4").ClearContents()

# Widen column B to fit the new, longer header text.
$ws.Columns.Item(2).ColumnWidth = 35

# Match the saved selection state from the authored workbook.
$ws.Range("B7").Select() | Out-Null
